$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Locate the two target paragraphs in the "Socks in the Dark" section
# robustly (by content), rather than assuming a fixed paragraph index.
# ------------------------------------------------------------------
$socksIdentifyIdx = -1
$socksEvaluateIdx = -1
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t.Contains("In the best- case scenario") -and $t.Contains("Identify potential solutions")) {
        $socksIdentifyIdx = $i
    }
    if ($t.StartsWith("Evaluate each potential solution:") -and $t.Length -lt 50) {
        $socksEvaluateIdx = $i
    }
}

if ($socksIdentifyIdx -eq -1) { throw "Could not locate the socks 'Identify potential solutions' paragraph" }
if ($socksEvaluateIdx -eq -1) { throw "Could not locate the socks 'Evaluate each potential solution' paragraph" }

# ------------------------------------------------------------------
# Paragraph: "Identify potential solutions:" (part b. sock-count math)
# Rewritten into several runs + keeps the existing _GoBack bookmark
# sitting between "Multiple pairs of " and "either the black...".
# ------------------------------------------------------------------
$pIdentify = $d.Paragraphs.Item($socksIdentifyIdx).Range
$identifyXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:mc="http://schemas.openxmlformats.org/markup-compatibility/2006" mc:Ignorable="w14"><w:body><w:p w14:paraId="542E8527" w14:textId="4842573F" w:rsidR="001979A6" w:rsidRDefault="001979A6" w:rsidP="001979A6"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="8"/></w:numPr><w:spacing w:line="480" w:lineRule="auto"/></w:pPr><w:r><w:t>Identify potential solutions:</w:t></w:r><w:r><w:br/><w:t xml:space="preserve">a. </w:t></w:r><w:r w:rsidR="00321500"><w:t>In the best-case scenario, a person would need 4 socks to come up with a guaranteed matching pair.  The first sock they drew plus the potential to draw one more of each color (black, brown or white) before they received a guaranteed pair.</w:t></w:r><w:r w:rsidR="00CC3668"><w:br/><w:t>b.  In the best- case scenario, a</w:t></w:r><w:r><w:t xml:space="preserve"> person would need a total of 18</w:t></w:r><w:r><w:t xml:space="preserve"> socks to come up with a guaranteed matching pair of socks in each color.  </w:t></w:r><w:r><w:t xml:space="preserve">Multiple pairs of </w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r w:rsidR="00CC3668"><w:t xml:space="preserve">either the black or brown could be drawn before a white pair was drawn.  Since there are only 4 white socks total and 2 of them need to be selected for a pair, all of the colored socks plus 2 of the white could potentially be selected prior to getting a matching white pair.  </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$null = $pIdentify.InsertXML($identifyXml)

# ------------------------------------------------------------------
# Paragraph: "Evaluate each potential solution:" (socks section)
# Rewritten into several runs, adding explanatory sentences to a/b.
# Re-locate by content since the InsertXML above may have shifted
# paragraph indices.
# ------------------------------------------------------------------
$socksEvaluateIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t.StartsWith("Evaluate each potential solution:") -and $t.Length -lt 50) {
        $socksEvaluateIdx = $i
    }
}
if ($socksEvaluateIdx -eq -1) { throw "Could not re-locate the socks 'Evaluate each potential solution' paragraph" }

$pEvaluate = $d.Paragraphs.Item($socksEvaluateIdx).Range
$evaluateXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:mc="http://schemas.openxmlformats.org/markup-compatibility/2006" mc:Ignorable="w14"><w:body><w:p w14:paraId="17B6F992" w14:textId="1CCECE89" w:rsidR="001979A6" w:rsidRDefault="001979A6" w:rsidP="001979A6"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="8"/></w:numPr><w:spacing w:line="480" w:lineRule="auto"/></w:pPr><w:r><w:t>Evalu</w:t></w:r><w:r><w:t>ate each potential solution:</w:t></w:r><w:r><w:br/><w:t xml:space="preserve">a.  Based on the mathematical probability, it would seem that each scenario would meet what the goal of the problem is.  </w:t></w:r><w:r w:rsidR="00A61AAB"><w:br/></w:r><w:r><w:t xml:space="preserve">b.  </w:t></w:r><w:r><w:t xml:space="preserve">The solutions are for the least amount of socks needed to make each statement’s guarantee.  Theoretically it seems that these solutions would work for all cases.  </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$null = $pEvaluate.InsertXML($evaluateXml)
